# Add six new playable races to the "Playable Races" table on Sheet1.
# Each new row is inserted immediately above the existing row for the
# race that should alphabetically follow it, then populated with the
# RACE (col A) / SUB-RACE (col B) / WORLD (col C) values.
#
# Step 1: do all the structural row inserts working from the BOTTOM of
# the sheet upward, so every insertion point still refers to the
# ORIGINAL (pre-edit) row numbers -- no offset bookkeeping needed.
#
# Step 2: fill in the new cells' values afterwards, in the same order
# the source workbook first introduces each new string (Fairy,
# Hobgoblin of the Feywild, Owlfolk, Rabbitfolk, Locathah, Verdan) so
# the shared-string table comes out laid out the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- structural inserts (bottom-up) ----------------------------------
$ws.Rows.Item(117).Insert()          # new row for Verdan, above old "Vedalkin"
$ws.Rows.Item(87).Insert()           # new row for Owlfolk, above old "Reborn"
$ws.Rows.Item(88).Insert()           # new row for Rabbitfolk, right after Owlfolk
$ws.Rows.Item(82).Insert()           # new row for Locathah, above old "Loxodon"
$ws.Rows.Item(74).Insert()           # new row for Hobgoblin of the Feywild, above old "Human"
$ws.Rows.Item(44).Insert()           # new row for Fairy, above old "Firbolg"

# -- populate values (in first-use order) -----------------------------
$ws.Cells.Item(44, 1).Value = "Fairy"
$ws.Cells.Item(44, 3).Value = "Feywild"

$ws.Cells.Item(75, 1).Value = "Hobgoblin of the Feywild"
$ws.Cells.Item(75, 3).Value = "Feywild"

$ws.Cells.Item(90, 1).Value = "Owlfolk"
$ws.Cells.Item(90, 3).Value = "Feywild"

$ws.Cells.Item(91, 1).Value = "Rabbitfolk"

$ws.Cells.Item(84, 1).Value = "Locathah"

$ws.Cells.Item(122, 1).Value = "Verdan"

# Update the view to match the edited workbook: scrolled down so row 91
# is the top visible row, with the active cell on A122 (the new Verdan
# row) selected.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 91
$aw.ScrollColumn = 1
[void]$ws.Range("A122").Select()
